$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Locate the three paragraphs involved in this edit by their text content
#    (robust to any index drift rather than relying on fixed paragraph numbers).
# ---------------------------------------------------------------------------
$licenseHeadingIdx = 0
$licenseTextIdx = 0
$pdfNoticeIdx = 0

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($licenseHeadingIdx -eq 0 -and $t -like "*License Information*") {
        $licenseHeadingIdx = $i
    }
    if ($licenseTextIdx -eq 0 -and $t -like "*is based on*") {
        $licenseTextIdx = $i
    }
    if ($pdfNoticeIdx -eq 0 -and $t -like "*This PDF version is provided under the same license*") {
        $pdfNoticeIdx = $i
    }
}

# ---------------------------------------------------------------------------
# 2. Remove the whole "License Information" Heading2 paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($licenseHeadingIdx).Range.Delete()

# Indices below $licenseHeadingIdx are unaffected; $licenseTextIdx and
# $pdfNoticeIdx both shift down by one paragraph.
$licenseTextIdx = $licenseTextIdx - 1
$pdfNoticeIdx = $pdfNoticeIdx - 1

# ---------------------------------------------------------------------------
# 3. Remove the "This PDF version is provided under the same license."
#    paragraph entirely - its content is being folded into the license
#    paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($pdfNoticeIdx).Range.Delete()

# ---------------------------------------------------------------------------
# 4. Replace the body of the license paragraph with the new wording. Leave
#    one character of "slack" at the tail so the paragraph's own trailing
#    (empty) run survives untouched, matching the original paragraph's
#    leading/trailing empty runs.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item($licenseTextIdx)
$clearRange = $d.Range($p.Range.Start, $p.Range.End - 1)
$clearRange.Delete()

$p = $d.Paragraphs.Item($licenseTextIdx)
$insertionPoint = $d.Range($p.Range.Start, $p.Range.Start)

$newRunsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t>unfoldingWord&#174; Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> &#169; 2022 unfoldingWord. Released under CC BY-SA 4.0 license. </w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t>unfoldingWord&#174; Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> has been adapted in the following languages: Tok Pisin, Arabic (&#1593;&#1585;&#1576;&#1610;), French (Fran&#231;ais), Hindi (&#2361;&#2367;&#2306;&#2342;&#2368;), Indonesian (Bahasa Indonesia), Portuguese (Portugu&#234;s), Russian (&#1056;&#1091;&#1089;&#1089;&#1082;&#1080;&#1081;), Spanish (Espa&#241;ol), Swahili (Kiswahili), and Simplified Chinese (&#31616;&#20307;&#20013;&#25991;) from </w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t>unfoldingWord&#174; Translation Questions</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t xml:space="preserve"> &#169; 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newRunsXml)

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
